$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 10,20

$data[0,0] = "Amanda"
$data[0,1] = "A"
$data[0,2] = "55-64"
$data[0,3] = "Masculino"
$data[0,4] = "Prefiro não dizer"
$data[0,5] = "Branca"
$data[0,6] = "Ateu"
$data[0,7] = "Centro Direita"
$data[0,8] = "Divorciado"
$data[0,9] = "Com pais ou responsáveis"
$data[0,10] = "1-3"
$data[0,11] = "9937"
$data[0,12] = "Cinema e Audiovisual"
$data[0,13] = "1"
$data[0,14] = "Não sei"
$data[0,15] = "Moderada"
$data[0,16] = "Alta"
$data[0,17] = "Baixa"
$data[0,18] = "Redes Sociais"
$data[0,19] = "Ausência das agências de checagem"

$data[1,0] = "Maria"
$data[1,1] = "C"
$data[1,2] = "65+"
$data[1,3] = "Masculino"
$data[1,4] = "Prefiro não dizer"
$data[1,5] = "Amarela"
$data[1,6] = "Evangélica"
$data[1,7] = "Esquerda"
$data[1,8] = "Viúvo"
$data[1,9] = "Outro"
$data[1,10] = "4-6"
$data[1,11] = "9994"
$data[1,12] = "Relações Públicas"
$data[1,13] = "6"
$data[1,14] = "Não sei"
$data[1,15] = "Baixa"
$data[1,16] = "Alta"
$data[1,17] = "Alta"
$data[1,18] = "Portais de notícias"
$data[1,19] = "A forma como a notícia é apresentada"

$data[2,0] = "Marcos"
$data[2,1] = "E"
$data[2,2] = "18-24"
$data[2,3] = "Masculino"
$data[2,4] = "Não binário"
$data[2,5] = "Branca"
$data[2,6] = "Espírita"
$data[2,7] = "Não se interessa"
$data[2,8] = "Solteiro"
$data[2,9] = "Sozinho"
$data[2,10] = "4-6"
$data[2,11] = "8271"
$data[2,12] = "Cinema e Audiovisual"
$data[2,13] = "1"
$data[2,14] = "Não sei"
$data[2,15] = "Moderada"
$data[2,16] = "Moderada"
$data[2,17] = "Moderada"
$data[2,18] = "Televisão e/ou Rádio"
$data[2,19] = "A forma como a notícia é apresentada"

$data[3,0] = "Thiago"
$data[3,1] = "E"
$data[3,2] = "25-34"
$data[3,3] = "Masculino"
$data[3,4] = "Não binário"
$data[3,5] = "Indígena"
$data[3,6] = "Católica"
$data[3,7] = "Centro Esquerda"
$data[3,8] = "Divorciado"
$data[3,9] = "Outro"
$data[3,10] = "7-10"
$data[3,11] = "8143"
$data[3,12] = "Publicidade e Propaganda"
$data[3,13] = "5"
$data[3,14] = "Não"
$data[3,15] = "Baixa"
$data[3,16] = "Alta"
$data[3,17] = "Baixa"
$data[3,18] = "Portais de notícias"
$data[3,19] = "Outro"

$data[4,0] = "Rayane"
$data[4,1] = "C"
$data[4,2] = "65+"
$data[4,3] = "Feminino"
$data[4,4] = "Outro"
$data[4,5] = "Amarela"
$data[4,6] = "Outras"
$data[4,7] = "Centro"
$data[4,8] = "Divorciado"
$data[4,9] = "Com parceiro e/ou filhos"
$data[4,10] = "<1"
$data[4,11] = "6765"
$data[4,12] = "Cinema e Audiovisual"
$data[4,13] = "4"
$data[4,14] = "Sim, muito diferentes"
$data[4,15] = "Alta"
$data[4,16] = "Alta"
$data[4,17] = "Moderada"
$data[4,18] = "Redes Sociais"
$data[4,19] = "As pessoas acreditam em notícias que lhes convêm"

$data[5,0] = "Ysadora"
$data[5,1] = "D"
$data[5,2] = "25-34"
$data[5,3] = "Feminino"
$data[5,4] = "Prefiro não dizer"
$data[5,5] = "Preta"
$data[5,6] = "Agnóstico"
$data[5,7] = "Não se interessa"
$data[5,8] = "Casado"
$data[5,9] = "Com parceiro e/ou filhos"
$data[5,10] = "1-3"
$data[5,11] = "6267"
$data[5,12] = "Publicidade e Propaganda"
$data[5,13] = "3"
$data[5,14] = "Sim, muito diferentes"
$data[5,15] = "Alta"
$data[5,16] = "Moderada"
$data[5,17] = "Alta"
$data[5,18] = "Aplicativos de mensagem"
$data[5,19] = "Polarização"

$data[6,0] = "Jonathan"
$data[6,1] = "D"
$data[6,2] = "<18"
$data[6,3] = "Feminino"
$data[6,4] = "Não binário"
$data[6,5] = "Prefiro não dizer"
$data[6,6] = "Católica"
$data[6,7] = "Direita"
$data[6,8] = "Divorciado"
$data[6,9] = "Sozinho"
$data[6,10] = "<1"
$data[6,11] = "7462"
$data[6,12] = "Relações Públicas"
$data[6,13] = "1"
$data[6,14] = "Não sei"
$data[6,15] = "Baixa"
$data[6,16] = "Moderada"
$data[6,17] = "Baixa"
$data[6,18] = "Aplicativos de mensagem"
$data[6,19] = "Compartilhamento por influenciadores"

$data[7,0] = "Talita"
$data[7,1] = "A"
$data[7,2] = "35-44"
$data[7,3] = "Masculino"
$data[7,4] = "Outro"
$data[7,5] = "Preta"
$data[7,6] = "Outras"
$data[7,7] = "Extrema Esquerda"
$data[7,8] = "Divorciado"
$data[7,9] = "Outro"
$data[7,10] = "1-3"
$data[7,11] = "5481"
$data[7,12] = "Jornalismo"
$data[7,13] = "8"
$data[7,14] = "Sim, porém pouco"
$data[7,15] = "Alta"
$data[7,16] = "Alta"
$data[7,17] = "Baixa"
$data[7,18] = "Redes Sociais"
$data[7,19] = "Polarização"

$data[8,0] = "Gabriel"
$data[8,1] = "D"
$data[8,2] = "18-24"
$data[8,3] = "Feminino"
$data[8,4] = "Não binário"
$data[8,5] = "Branca"
$data[8,6] = "Judaíca"
$data[8,7] = "Centro Direita"
$data[8,8] = "Prefiro não dizer"
$data[8,9] = "Com parceiro e/ou filhos"
$data[8,10] = ">10"
$data[8,11] = "5050"
$data[8,12] = "Cinema e Audiovisual"
$data[8,13] = "8"
$data[8,14] = "Não são diferentes"
$data[8,15] = "Baixa"
$data[8,16] = "Alta"
$data[8,17] = "Alta"
$data[8,18] = "Televisão e/ou Rádio"
$data[8,19] = "Compartilhamento por influenciadores"

$data[9,0] = "Gabriel"
$data[9,1] = "D"
$data[9,2] = "45-54"
$data[9,3] = "Masculino"
$data[9,4] = "Transgênero"
$data[9,5] = "Branca"
$data[9,6] = "Ateu"
$data[9,7] = "Extrema Esquerda"
$data[9,8] = "Divorciado"
$data[9,9] = "Outro"
$data[9,10] = ">10"
$data[9,11] = "7461"
$data[9,12] = "Relações Públicas"
$data[9,13] = "2"
$data[9,14] = "Sim, porém pouco"
$data[9,15] = "Baixa"
$data[9,16] = "Baixa"
$data[9,17] = "Baixa"
$data[9,18] = "Aplicativos de mensagem"
$data[9,19] = "Outro"

$ws.Range("B2:U11").Value = $data